# Corrige/actualiza la base de datos de Estado de Cuenta:
# intercambia los periodos de mora (y sus valores) entre las filas 16 y 17.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fila 16: Periodo Mora 2402 / Valor Mora 34666  ->  2401 / 3467
# Fila 17: Periodo Mora 2401 / Valor Mora 3467    ->  2402 / 34666
$ws.Range("E16").Value = "2401"
$ws.Range("F16").Value = 3467

$ws.Range("E17").Value = "2402"
$ws.Range("F17").Value = 34666
